# Add the new "Pentester Lab" worksheet immediately before "Udemy"
# (Worksheets.Add inserts before the sheet passed as argument and makes
# the new sheet the active one, matching the diff's activeTab bump and
# the tabSelected flag moving off "Udacity" onto the new sheet.)
$wb = $excel.ActiveWorkbook
$udemySheet = $wb.Worksheets.Item("Udemy")
$hackerOneSheet = $wb.Worksheets.Item("Hacker One")

$psheet = $wb.Worksheets.Add($udemySheet)
$psheet.Name = "Pentester Lab"

# --- content -------------------------------------------------------
$psheet.Range("B2").Value = "Android Badge"
$psheet.Range("C2").Value = "https://pentesterlab.com/badges/android"

$psheet.Range("B4").Value = "HTTP Badge"
$psheet.Range("C4").Value = "https://pentesterlab.com/badges/http"

$psheet.Range("B6").Value = "API Badge"
$psheet.Range("C6").Value = "https://pentesterlab.com/badges/api"

$psheet.Range("B8").Value = "Essential Badge"
$psheet.Range("C8").Value = "https://pentesterlab.com/badges/essential"

# --- formatting ------------------------------------------------------
# "Hacker One"!B2 already carries the Segoe UI font used for the badge
# labels; copy its format onto B2 first, then finish resolving the
# alignment (left / center / wrap) on that single cell so only one new
# cellXf is produced, and finally fan that fully-resolved format out to
# the other three label cells with a single paste each.
[void]$hackerOneSheet.Range("B2").Copy()
[void]$psheet.Range("B2").PasteSpecial(-4122)
$psheet.Range("B2").HorizontalAlignment = -4131
$psheet.Range("B2").VerticalAlignment = -4108
$psheet.Range("B2").WrapText = $true

[void]$psheet.Range("B2").Copy()
[void]$psheet.Range("B4").PasteSpecial(-4122)
[void]$psheet.Range("B6").PasteSpecial(-4122)
[void]$psheet.Range("B8").PasteSpecial(-4122)

# --- column widths / row heights -------------------------------------
$psheet.Columns.Item(2).ColumnWidth = 52.5
$psheet.Columns.Item(3).ColumnWidth = 70

$psheet.Rows.Item(2).RowHeight = 16.8
$psheet.Rows.Item(4).RowHeight = 16.8
$psheet.Rows.Item(6).RowHeight = 16.8
$psheet.Rows.Item(8).RowHeight = 16.8

# --- selection on the new sheet ---------------------------------------
$psheet.Range("C8").Select()
